# Generate Report for Handback
# Update "Latest Handback DateTime" (column K) for the row-2 file
# (1169e9a6-ecc9-4d69-a72d-10aa163b9c7a) on both the zh-cn and de-de
# localization-status sheets to reflect the new handback timestamps.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("K2").Value = "2016-09-06 15:51:47"

$de = $wb.Worksheets.Item("de-de")
$de.Range("K2").Value = "2016-09-06 15:52:15"
